# Generate Report for Handback
# Updates the localization-status workbook to reflect that file "a.md" has
# been handed back (target file / handback xlf / handback datetime populated)
# and the Status changes from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is shown (Overview zh-cn/
# de-de columns as well as the per-language "Status" column).

$wb = $excel.ActiveWorkbook

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8db7f4ba2bf9ba57cdd2327cc6107eefe289827f/e2e/a.md"

# Hyperlink-like font formatting (matches the workbook's existing custom
# "HyperLink" cell style: underlined, Cornflower Blue FF6495ED).
$hyperlinkColor = 15570276   # BGR encoding of RGB(0x64,0x95,0xED)

# ---------------------------------------------------------------------
# 1) Overview sheet: Status columns (zh-cn / de-de) for both rows move
#    from "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen columns E & F so the longer status text fits (was ~17.2 chars).
$overview.Range("E1").ColumnWidth = 29.17
$overview.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# 2) zh-cn sheet: row for a.md (row 2) and row for b.md (row 3) both get
#    Latest Target File / Latest Handback File / Latest Handback DateTime
#    populated (the report tool stamps both rows once a.md is handed back).
#    The per-language "Status" column (C) mirrors the Overview status text,
#    so it also moves from "Ready for handoff" to the handed-back wording.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aUrl, "", "", "a.md")
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-20 18:48:03"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aUrl, "", "", "a.md")
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-20 18:48:03"

# Widen Status (C) and Latest Handback File (J) columns.
$zhcn.Range("C1").ColumnWidth = 29.17
$zhcn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3) de-de sheet: same shape of update as zh-cn, different handback xlf /
#    datetime values.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Hyperlinks.Add($dede.Range("I2"), $aUrl, "", "", "a.md")
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-20 18:48:11"

$dede.Hyperlinks.Add($dede.Range("I3"), $aUrl, "", "", "a.md")
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-20 18:48:11"

# Widen Status (C) and Latest Handback File (J) columns.
$dede.Range("C1").ColumnWidth = 29.17
$dede.Range("J1").ColumnWidth = 39.17
